$d = $word.ActiveDocument

# --- Step 1: merge "Description entrées/sorties GP" + "IO" into a single
# run reading "Description entrées/sorties GPIO", dropping the old
# _GoBack bookmark that used to sit between the two runs.
# (Overwriting a range with text that reads identically to what is
# already there is treated as a no-op by the engine, so first overwrite
# with a throwaway placeholder to force the runs to actually be
# rewritten/merged, then set the real final text.)
$p4 = $d.Paragraphs(4)
$rng4 = $d.Range($p4.Range.Start, $p4.Range.End - 1)
$rng4.Text = "__PLACEHOLDER__"

$p4 = $d.Paragraphs(4)
$rng4 = $d.Range($p4.Range.Start, $p4.Range.End - 1)
$rng4.Text = "Description entrées/sorties GPIO"

# --- Step 2: add the new bullet paragraph right after it, same list style.
$p4 = $d.Paragraphs(4)
$p4.Range.InsertParagraphAfter()

$p5 = $d.Paragraphs(5)
$rng5 = $d.Range($p5.Range.Start, $p5.Range.End - 1)
# A trailing placeholder char is appended for now; it lets us anchor the
# bookmark right after the real text without hitting the "end of
# document" collapsed-range edge case (see step 3), and gets stripped
# back off afterwards.
$rng5.Text = "Convertisseur analogique numérique~"

# --- Step 3: put the _GoBack bookmark right after that text, before the
# trailing placeholder char (there is only ever one _GoBack bookmark in
# a document, so adding it here moves it away from paragraph 4).
$p5 = $d.Paragraphs(5)
$bmPos = $p5.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- Step 4: strip the trailing placeholder char back off.
$p5 = $d.Paragraphs(5)
$trailRange = $d.Range($p5.Range.End - 2, $p5.Range.End - 1)
$trailRange.Delete()
